$d = $word.ActiveDocument
$t = $d.Tables(1)

# Row 2 (1-indexed): add center alignment to existing content, no text changes
$t.Cell(2,1).Range.Paragraphs(1).Alignment = 1
$t.Cell(2,2).Range.Paragraphs(1).Alignment = 1
$t.Cell(2,3).Range.Paragraphs(1).Alignment = 1
$t.Cell(2,4).Range.Paragraphs(1).Alignment = 1

# Row 3 (1-indexed): add center alignment + new text content
$cellDatum = $t.Cell(3,1)
$cellDatum.Range.Paragraphs(1).Alignment = 1
$rngDatum = $cellDatum.Range
$rngDatum.Text = "5.6.2021"
$rngDatumFmt = $cellDatum.Range
$rngDatumFmt.Font.Size = 12

$cellVerzija = $t.Cell(3,2)
$cellVerzija.Range.Paragraphs(1).Alignment = 1
$rngVerzija = $cellVerzija.Range
$rngVerzija.Text = "1.0"
$rngVerzijaFmt = $cellVerzija.Range
$rngVerzijaFmt.Font.Size = 12

$cellOpis = $t.Cell(3,3)
$cellOpis.Range.Paragraphs(1).Alignment = 1
$rngOpis = $cellOpis.Range
$rngOpis.Text = "Finalna verzija"
$rngOpisFmt = $cellOpis.Range
$rngOpisFmt.Font.Size = 12

# Autor cell needs multiple runs (name, break+surname, separate last letter with language tag)
$cellAutor = $t.Cell(3,4)
$pAutor = $cellAutor.Range.Paragraphs(1)
$rngAutor = $pAutor.Range
$xmlAutor = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="08912CCF" w14:textId="77777777" w:rsidR="0067036C" w:rsidRDefault="0067036C"><w:pPr><w:pStyle w:val="TableParagraph"/><w:jc w:val="center"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Marko Gloginja</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:br/><w:t>Stefan Lukovi</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>ć</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rngAutor.InsertXML($xmlAutor) | Out-Null
